# Add two more demo rows to the Product_Name sheet (Ate/Atenolol and
# Colchi/Colchisin), matching the style/format of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (cell styles, number format) of the last existing
# data row (row 3) down onto the two new rows (4 and 5) before filling in
# values, so the new cells pick up the same styles (bold-free text style /
# date number format) as the rest of the table.
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C5").PasteSpecial(-4122)

$ws.Range("A4").Value = "Ate"
$ws.Range("B4").Value = "Atenolol"
$ws.Range("C4").Value = 44227

$ws.Range("A5").Value = "Colchi"
$ws.Range("B5").Value = "Colchisin"
$ws.Range("C5").Value = 44227
